# Insert a new data row right above current row 322 (Ají price record),
# shifting the existing rows 322-415 down to 323-416, and fill the new
# row 322 with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 322; everything below shifts down by one.
$ws.Rows.Item(322).Insert()

# Populate the newly inserted row 322 with the new record's data.
$ws.Range("A322").Value = 3
$ws.Range("B322").Value = 'Femacal de La Calera'
$ws.Range("C322").Value = 'Coquimbo'
$ws.Range("D322").Value = 44642
$ws.Range("E322").Value = 5
$ws.Range("F322").Value = 100112021
$ws.Range("G322").Value = 'Ají'
$ws.Range("H322").Value = 'Americana (o)'
$ws.Range("I322").Value = 'Primera'
$ws.Range("J322").Value = 76
$ws.Range("K322").Value = 14000
$ws.Range("L322").Value = 14500
$ws.Range("M322").Value = 14250
$ws.Range("N322").Value = '$/caja 15 kilos'
$ws.Range("O322").Value = 'Limache'
$ws.Range("P322").Value = 950
$ws.Range("Q322").Value = 15
$ws.Range("R322").Value = 'Hortaliza'
